$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new price record was inserted before the current row 95, shifting all
# subsequent rows (95-228) down by one (the former last row, 228, becomes 229).
$ws.Rows.Item(95).Insert()

# Populate the newly inserted row 95 with the new weekly record.
$ws.Range("A95").Value = 8
$ws.Range("B95").Value = "Terminal La Palmera de La Serena"
$ws.Range("C95").Value = "Coquimbo"
$ws.Range("D95").Value = 44792
$ws.Range("E95").Value = 4
$ws.Range("F95").Value = 100112037
$ws.Range("G95").Value = "Cebollín"
$ws.Range("H95").Value = "Sin especificar"
$ws.Range("I95").Value = "Primera"
$ws.Range("J95").Value = 2000
$ws.Range("K95").Value = 1400
$ws.Range("L95").Value = 1600
$ws.Range("M95").Value = 1500
$ws.Range("N95").Value = "$/paquete 6 unidades"
$ws.Range("O95").Value = "Provincia del Elquí"
$ws.Range("P95").Value = 250
$ws.Range("Q95").Value = 6
$ws.Range("R95").Value = "Hortaliza"
